$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Updates" change-log entries in column A for rows 4 through 13
# (originele situatie, tuples ipv lijsten, efficiente stappen, ...,
# vermijden occupiedBy()) are obsolete now that the heuristic work has
# moved on to the Heapq-based approach, so clear them out.
$ws.Range("A4:A13").ClearContents()

# Move / leave the selection where the edit left off.
$ws.Range("A14").Select()
